$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds rows 168-180 of "Fruta / hortaliza" price data.
# This weekly update inserts two new observation rows at the top of that
# block (new row 168 and 169), pushing the existing rows 168-180 down to
# 170-182.
$ws.Rows.Item(168).Insert()
$ws.Rows.Item(168).Insert()

# --- New row 168 ---
$ws.Cells.Item(168, 1).Value = 6
$ws.Cells.Item(168, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(168, 3).Value = "Metropolitana"
$ws.Cells.Item(168, 4).Value = 44491
$ws.Cells.Item(168, 5).Value = 13
$ws.Cells.Item(168, 6).Value = "Fruta"
$ws.Cells.Item(168, 7).Value = 100101
$ws.Cells.Item(168, 8).Value = "Berries"
$ws.Cells.Item(168, 9).Value = 100101001
$ws.Cells.Item(168, 10).Value = "Arándano (blue)"
$ws.Cells.Item(168, 11).Value = "Sin especificar"
$ws.Cells.Item(168, 12).Value = "Primera"
$ws.Cells.Item(168, 13).Value = 850
$ws.Cells.Item(168, 14).Value = 12000
$ws.Cells.Item(168, 15).Value = 12000
$ws.Cells.Item(168, 16).Value = 12000
$ws.Cells.Item(168, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(168, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 19).Value = 6000
$ws.Cells.Item(168, 20).Value = 2

# --- New row 169 ---
$ws.Cells.Item(169, 1).Value = 6
$ws.Cells.Item(169, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(169, 3).Value = "Metropolitana"
$ws.Cells.Item(169, 4).Value = 44491
$ws.Cells.Item(169, 5).Value = 13
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100101
$ws.Cells.Item(169, 8).Value = "Berries"
$ws.Cells.Item(169, 9).Value = 100101001
$ws.Cells.Item(169, 10).Value = "Arándano (blue)"
$ws.Cells.Item(169, 11).Value = "Sin especificar"
$ws.Cells.Item(169, 12).Value = "Segunda"
$ws.Cells.Item(169, 13).Value = 150
$ws.Cells.Item(169, 14).Value = 10000
$ws.Cells.Item(169, 15).Value = 10000
$ws.Cells.Item(169, 16).Value = 10000
$ws.Cells.Item(169, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(169, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(169, 19).Value = 5000
$ws.Cells.Item(169, 20).Value = 2
